# Auto-generated script to update Goblin Profits market data cells
# Applies cached market-price / profit value refresh across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 188.4
$ws.Range("I4").Value = 188.4
$ws.Range("K4").Value = 188.4
$ws.Range("M4").Value = -74.40000000000001
$ws.Range("H6").Value = 5766.5557
$ws.Range("I6").Value = 6448.6875
$ws.Range("J6").Value = 309.5
$ws.Range("K6").Value = 19346.0625
$ws.Range("L6").Value = 928.5
$ws.Range("M6").Value = -19234.0625
$ws.Range("N6").Value = -1152.5
$ws.Range("H12").Value = 641
$ws.Range("I12").Value = 561
$ws.Range("J12").Value = 721
$ws.Range("K12").Value = 561
$ws.Range("L12").Value = 721
$ws.Range("M12").Value = -391
$ws.Range("N12").Value = -1061
$ws.Range("H28").Value = 737.4400000000001
$ws.Range("I28").Value = 734
$ws.Range("K28").Value = 734
$ws.Range("M28").Value = -249
$ws.Range("H40").Value = 2825.7778
$ws.Range("I40").Value = 1220.3
$ws.Range("J40").Value = 3770.1765
$ws.Range("K40").Value = 1220.3
$ws.Range("L40").Value = 3770.1765
$ws.Range("M40").Value = -1045.3
$ws.Range("N40").Value = -4120.1765
$ws.Range("H41").Value = 735.625
$ws.Range("I41").Value = 735.625
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 735.625
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -295.625
$ws.Range("N41").ClearContents()
$ws.Range("H88").Value = 6201.4
$ws.Range("I88").Value = 2253.6667
$ws.Range("J88").Value = 8833.223
$ws.Range("K88").Value = 2253.6667
$ws.Range("L88").Value = 8833.223
$ws.Range("M88").Value = -1847.6667
$ws.Range("N88").Value = -9645.223
$ws.Range("H91").Value = 6201.4
$ws.Range("I91").Value = 2253.6667
$ws.Range("J91").Value = 8833.223
$ws.Range("K91").Value = 2253.6667
$ws.Range("L91").Value = 8833.223
$ws.Range("M91").Value = -849.6667000000002
$ws.Range("N91").Value = -11641.223
$ws.Range("H92").Value = 1999.2941
$ws.Range("I92").Value = 2082
$ws.Range("K92").Value = 2082
$ws.Range("M92").Value = -834
$ws.Range("H107").Value = 964.89655
$ws.Range("I107").Value = 1054.7391
$ws.Range("K107").Value = 1054.7391
$ws.Range("M107").Value = 865.2609
$ws.Range("H137").Value = 660387.9399999999
$ws.Range("I137").Value = 1267.2667
$ws.Range("K137").Value = 3801.800099999999
$ws.Range("M137").Value = -1251.800099999999
$ws.Range("H138").Value = 1955594.6
$ws.Range("I138").Value = 1463.5217
$ws.Range("J138").Value = 3277506.8
$ws.Range("K138").Value = 4390.5651
$ws.Range("L138").Value = 9832520.399999999
$ws.Range("M138").Value = 749.4349000000002
$ws.Range("N138").Value = -9842800.399999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1159.4584
$ws.Range("I2").Value = 379.44446
$ws.Range("K2").Value = 379.44446
$ws.Range("M2").Value = -266.44446
$ws.Range("H13").Value = 67497.664
$ws.Range("J13").Value = 1246.5
$ws.Range("L13").Value = 1246.5
$ws.Range("N13").Value = -1534.5
$ws.Range("H25").Value = 1070.5
$ws.Range("I25").Value = 1070.5
$ws.Range("K25").Value = 1070.5
$ws.Range("M25").Value = -668.5
$ws.Range("H32").Value = 79585.88
$ws.Range("I32").Value = 80640.21000000001
$ws.Range("K32").Value = 80640.21000000001
$ws.Range("M32").Value = -80353.21000000001
$ws.Range("H63").Value = 7490.6875
$ws.Range("I63").Value = 4550.1113
$ws.Range("K63").Value = 4550.1113
$ws.Range("M63").Value = -3864.1113
$ws.Range("H66").Value = 7490.6875
$ws.Range("I66").Value = 4550.1113
$ws.Range("K66").Value = 22750.5565
$ws.Range("M66").Value = -19318.5565
$ws.Range("H97").Value = 385.625
$ws.Range("I97").Value = 375.83334
$ws.Range("K97").Value = 375.83334
$ws.Range("M97").Value = 120.16666
$ws.Range("H116").Value = 1159.4584
$ws.Range("I116").Value = 379.44446
$ws.Range("K116").Value = 379.44446
$ws.Range("M116").Value = 1914.55554

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1159.4584
$ws.Range("I3").Value = 379.44446
$ws.Range("K3").Value = 379.44446
$ws.Range("M3").Value = -265.44446
$ws.Range("H94").Value = 1960.9231
$ws.Range("I94").Value = 2199.4546
$ws.Range("J94").Value = 649
$ws.Range("K94").Value = 2199.4546
$ws.Range("L94").Value = 649
$ws.Range("M94").Value = -1748.4546
$ws.Range("N94").Value = -1551
$ws.Range("H107").Value = 4623.0293
$ws.Range("J107").Value = 5913.2
$ws.Range("L107").Value = 5913.2
$ws.Range("N107").Value = -9753.200000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1807.3414
$ws.Range("I58").Value = 1785.5143
$ws.Range("K58").Value = 1785.5143
$ws.Range("M58").Value = -1582.5143
$ws.Range("H132").Value = 4495.5186
$ws.Range("I132").Value = 1518.9
$ws.Range("K132").Value = 4556.700000000001
$ws.Range("M132").Value = -2026.700000000001
$ws.Range("H136").Value = 1807.3414
$ws.Range("I136").Value = 1785.5143
$ws.Range("K136").Value = 5356.5429
$ws.Range("M136").Value = -2806.5429

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 92.30768999999999
$ws.Range("I14").Value = 92.30768999999999
$ws.Range("K14").Value = 276.92307
$ws.Range("M14").Value = -103.92307
$ws.Range("H17").Value = 5556471.5
$ws.Range("J17").Value = 1222.8334
$ws.Range("L17").Value = 3668.5002
$ws.Range("N17").Value = -4006.5002
$ws.Range("H80").Value = 3624.75
$ws.Range("J80").Value = 3500
$ws.Range("L80").Value = 10500
$ws.Range("N80").Value = -12372
$ws.Range("H83").Value = 3624.75
$ws.Range("J83").Value = 3500
$ws.Range("L83").Value = 31500
$ws.Range("N83").Value = -40860
$ws.Range("H86").Value = 780.6429000000001
$ws.Range("I86").Value = 629
$ws.Range("J86").Value = 1159.75
$ws.Range("K86").Value = 1887
$ws.Range("L86").Value = 3479.25
$ws.Range("M86").Value = -701
$ws.Range("N86").Value = -5851.25
$ws.Range("H89").Value = 780.6429000000001
$ws.Range("I89").Value = 629
$ws.Range("J89").Value = 1159.75
$ws.Range("K89").Value = 5661
$ws.Range("L89").Value = 10437.75
$ws.Range("M89").Value = 267
$ws.Range("N89").Value = -22293.75
$ws.Range("H113").Value = 1628
$ws.Range("J113").Value = 1222.2222
$ws.Range("L113").Value = 3666.6666
$ws.Range("N113").Value = -8006.6666
$ws.Range("H114").Value = 1902.2222
$ws.Range("J114").Value = 3745
$ws.Range("L114").Value = 11235
$ws.Range("N114").Value = -17743

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 361.25
$ws.Range("I9").Value = 126.42857
$ws.Range("J9").Value = 690
$ws.Range("K9").Value = 126.42857
$ws.Range("L9").Value = 690
$ws.Range("M9").Value = 43.57143000000001
$ws.Range("N9").Value = -1030
$ws.Range("H63").Value = 750007500
$ws.Range("J63").Value = 750007500
$ws.Range("L63").Value = 750007500
$ws.Range("N63").Value = -750008872
$ws.Range("H66").Value = 750007500
$ws.Range("J66").Value = 750007500
$ws.Range("L66").Value = 2250022500
$ws.Range("N66").Value = -2250029364
$ws.Range("H132").Value = 32259436
$ws.Range("I132").Value = 37038076
$ws.Range("K132").Value = 111114228
$ws.Range("M132").Value = -111111698

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 372499.16
$ws.Range("I2").Value = 372499.16
$ws.Range("K2").Value = 372499.16
$ws.Range("M2").Value = -372387.16
$ws.Range("H3").Value = 3000
$ws.Range("I3").Value = 3000
$ws.Range("K3").Value = 3000
$ws.Range("M3").Value = -2888
$ws.Range("H12").Value = 4499.5
$ws.Range("J12").Value = 4499.5
$ws.Range("L12").Value = 4499.5
$ws.Range("N12").Value = -4839.5
$ws.Range("H15").Value = 3000
$ws.Range("I15").Value = 3000
$ws.Range("K15").Value = 3000
$ws.Range("M15").Value = -2830
$ws.Range("H40").Value = 4511.7666
$ws.Range("I40").Value = 3430.7778
$ws.Range("J40").Value = 6133.25
$ws.Range("K40").Value = 3430.7778
$ws.Range("L40").Value = 6133.25
$ws.Range("M40").Value = -3294.7778
$ws.Range("N40").Value = -6405.25
$ws.Range("H46").Value = 851.7
$ws.Range("I46").Value = 440.125
$ws.Range("K46").Value = 440.125
$ws.Range("M46").Value = -252.125
$ws.Range("H93").Value = 4621.3784
$ws.Range("I93").Value = 2778.5264
$ws.Range("J93").Value = 6566.6113
$ws.Range("K93").Value = 2778.5264
$ws.Range("L93").Value = 6566.6113
$ws.Range("M93").Value = -1530.5264
$ws.Range("N93").Value = -9062.6113
$ws.Range("H136").Value = 44022.934
$ws.Range("I136").Value = 5116.5557
$ws.Range("J136").Value = 102382.5
$ws.Range("K136").Value = 15349.6671
$ws.Range("L136").Value = 307147.5
$ws.Range("M136").Value = -12799.6671
$ws.Range("N136").Value = -312247.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1692.4615
$ws.Range("I81").Value = 1780.2
$ws.Range("J81").Value = 1400
$ws.Range("K81").Value = 3560.4
$ws.Range("L81").Value = 2800
$ws.Range("M81").Value = -2499.4
$ws.Range("N81").Value = -4922
$ws.Range("H84").Value = 1692.4615
$ws.Range("I84").Value = 1780.2
$ws.Range("J84").Value = 1400
$ws.Range("K84").Value = 17802
$ws.Range("L84").Value = 14000
$ws.Range("M84").Value = -12498
$ws.Range("N84").Value = -24608
$ws.Range("H94").Value = 78191.5
$ws.Range("J94").Value = 78191.5
$ws.Range("L94").Value = 78191.5
$ws.Range("N94").Value = -79993.5
$ws.Range("H107").Value = 29412448
$ws.Range("I107").Value = 719.25
$ws.Range("J107").Value = 55556210
$ws.Range("K107").Value = 2157.75
$ws.Range("L107").Value = 166668630
$ws.Range("M107").Value = -237.75
$ws.Range("N107").Value = -166672470
